$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# ---- Add new data rows (45-60) to the Data-Org sheet ----

# Row 45
$ws.Range("A45").Value = 3368
$ws.Range("B45").Value = 4
$ws.Range("C45").Value = 5
$ws.Range("D45").Value = 3
$ws.Range("E45").Value = 5
$ws.Range("F45").Value = 5
$ws.Range("G45").Value = 4
$ws.Range("H45").Value = 2
$ws.Range("I45").Value = 2
$ws.Range("J45").Value = 4
$ws.Range("K45").Value = 5
$ws.Range("L45").Value = 4
$ws.Range("M45").Value = 4
$ws.Range("N45").Value = 4
$ws.Range("O45").Value = 5
$ws.Range("R45").Value = 3
$ws.Range("T45").Value = 3
$ws.Range("W45").Value = 4
$ws.Range("X45").Value = 5
$ws.Range("Y45").Value = 3
$ws.Range("AC45").Value = 5
$ws.Range("AD45").Value = 4
$ws.Range("AE45").Value = 5
$ws.Range("AF45").Value = 2
$ws.Range("AH45").Value = 4
$ws.Range("AJ45").Value = 1
$ws.Range("AL45").Value = 3
$ws.Range("A1").Copy()
$ws.Range("A45:B45").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D45:P45").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("R45").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("X45:Y45").PasteSpecial(-4122)

# Row 46
$ws.Range("A46").Value = 2456
$ws.Range("B46").Value = 5
$ws.Range("C46").Value = 5
$ws.Range("D46").Value = 2
$ws.Range("G46").Value = 2
$ws.Range("I46").Value = 5
$ws.Range("K46").Value = 5
$ws.Range("L46").Value = 4
$ws.Range("M46").Value = 4
$ws.Range("O46").Value = 5
$ws.Range("Q46").Value = 3
$ws.Range("S46").Value = 1
$ws.Range("T46").Value = 4
$ws.Range("U46").Value = 4
$ws.Range("V46").Value = 3
$ws.Range("W46").Value = 4
$ws.Range("AA46").Value = 3
$ws.Range("AC46").Value = 5
$ws.Range("AF46").Value = 3
$ws.Range("AJ46").Value = 4
$ws.Range("AL46").Value = 3
$ws.Range("A1").Copy()
$ws.Range("A46:B46").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G46").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I46").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("M46").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("O46").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("Q46").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("U46").PasteSpecial(-4122)

# Row 47
$ws.Range("A47").Value = 3436
$ws.Range("B47").Value = 5
$ws.Range("C47").Value = 4
$ws.Range("D47").Value = 5
$ws.Range("E47").Value = 5
$ws.Range("F47").Value = 3
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 3
$ws.Range("I47").Value = 5
$ws.Range("J47").Value = 4
$ws.Range("K47").Value = 4
$ws.Range("L47").Value = 3
$ws.Range("O47").Value = 4
$ws.Range("P47").Value = 4
$ws.Range("Q47").Value = 5
$ws.Range("R47").Value = 1
$ws.Range("T47").Value = 3
$ws.Range("V47").Value = 4
$ws.Range("X47").Value = 4
$ws.Range("Y47").Value = 4
$ws.Range("Z47").Value = 3
$ws.Range("AB47").Value = 4
$ws.Range("AC47").Value = 5
$ws.Range("AE47").Value = 3
$ws.Range("AI47").Value = 3
$ws.Range("AJ47").Value = 4
$ws.Range("AK47").Value = 5
$ws.Range("A1").Copy()
$ws.Range("A47:B47").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D47:G47").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I47:J47").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("O47").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("X47:Y47").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("AB47").PasteSpecial(-4122)

# Row 48
$ws.Range("A48").Value = 2356
$ws.Range("B48").Value = 4
$ws.Range("C48").Value = 5
$ws.Range("D48").Value = 3
$ws.Range("G48").Value = 2
$ws.Range("I48").Value = 4
$ws.Range("K48").Value = 4
$ws.Range("M48").Value = 4
$ws.Range("N48").Value = 4
$ws.Range("O48").Value = 3
$ws.Range("P48").Value = 2
$ws.Range("U48").Value = 4
$ws.Range("W48").Value = 5
$ws.Range("Y48").Value = 3
$ws.Range("Z48").Value = 2
$ws.Range("AC48").Value = 2
$ws.Range("AE48").Value = 3
$ws.Range("AH48").Value = 4
$ws.Range("AJ48").Value = 5
$ws.Range("AK48").Value = 1
$ws.Range("AL48").Value = 3
$ws.Range("A1").Copy()
$ws.Range("A48:B48").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G48").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I48").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("M48").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("O48").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("U48").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("Y48").PasteSpecial(-4122)

# Row 49
$ws.Range("A49").Value = 8903
$ws.Range("B49").Value = 5
$ws.Range("E49").Value = 5
$ws.Range("F49").Value = 2
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 4
$ws.Range("I49").Value = 5
$ws.Range("J49").Value = 5
$ws.Range("L49").Value = 4
$ws.Range("N49").Value = 4
$ws.Range("O49").Value = 5
$ws.Range("P49").Value = 5
$ws.Range("Q49").Value = 4
$ws.Range("R49").Value = 1
$ws.Range("S49").Value = 4
$ws.Range("T49").Value = 5
$ws.Range("U49").Value = 5
$ws.Range("V49").Value = 4
$ws.Range("X49").Value = 4
$ws.Range("AA49").Value = 4
$ws.Range("AE49").Value = 3
$ws.Range("AF49").Value = 5
$ws.Range("AG49").Value = 4
$ws.Range("AH49").Value = 2
$ws.Range("AJ49").Value = 3
$ws.Range("AL49").Value = 4
$ws.Range("A1").Copy()
$ws.Range("A49:B49").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("E49:G49").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I49:J49").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("O49").PasteSpecial(-4122)

# Row 50
$ws.Range("A50").Value = 2089
$ws.Range("B50").Value = 3
$ws.Range("D50").Value = 5
$ws.Range("E50").Value = 4
$ws.Range("G50").Value = 4
$ws.Range("I50").Value = 4
$ws.Range("O50").Value = 5
$ws.Range("Q50").Value = 4
$ws.Range("S50").Value = 3
$ws.Range("T50").Value = 5
$ws.Range("Y50").Value = 2
$ws.Range("AA50").Value = 3
$ws.Range("AC50").Value = 5
$ws.Range("AD50").Value = 5
$ws.Range("AE50").Value = 4
$ws.Range("AF50").Value = 5
$ws.Range("AG50").Value = 5
$ws.Range("AH50").Value = 3
$ws.Range("AI50").Value = 2
$ws.Range("AJ50").Value = 3
$ws.Range("AK50").Value = 1
$ws.Range("AL50").Value = 4
$ws.Range("A1").Copy()
$ws.Range("A50:B50").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("C50").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D50:E50").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G50").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("H50").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I50").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("J50:N50").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("O50").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("Y50").PasteSpecial(-4122)

# Row 51
$ws.Range("A51").Value = 7258
$ws.Range("B51").Value = 4
$ws.Range("C51").Value = 5
$ws.Range("E51").Value = 5
$ws.Range("I51").Value = 4
$ws.Range("J51").Value = 4
$ws.Range("K51").Value = 4
$ws.Range("L51").Value = 3
$ws.Range("M51").Value = 3
$ws.Range("N51").Value = 3
$ws.Range("P51").Value = 4
$ws.Range("R51").Value = 1
$ws.Range("T51").Value = 3
$ws.Range("V51").Value = 4
$ws.Range("X51").Value = 3
$ws.Range("Z51").Value = 4
$ws.Range("AA51").Value = 4
$ws.Range("AB51").Value = 4
$ws.Range("AC51").Value = 5
$ws.Range("AE51").Value = 5
$ws.Range("AF51").Value = 5
$ws.Range("AI51").Value = 4
$ws.Range("AK51").Value = 3
$ws.Range("A1").Copy()
$ws.Range("A51:B51").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("E51").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I51").PasteSpecial(-4122)

# Row 52
$ws.Range("A52").Value = 1299
$ws.Range("B52").Value = 4
$ws.Range("C52").Value = 5
$ws.Range("D52").Value = 3
$ws.Range("E52").Value = 5
$ws.Range("F52").Value = 5
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 1
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 5
$ws.Range("K52").Value = 5
$ws.Range("L52").Value = 3
$ws.Range("M52").Value = 4
$ws.Range("N52").Value = 5
$ws.Range("O52").Value = 5
$ws.Range("P52").Value = 5
$ws.Range("Q52").Value = 5
$ws.Range("R52").Value = 4
$ws.Range("S52").Value = 2
$ws.Range("T52").Value = 2
$ws.Range("U52").Value = 1
$ws.Range("V52").Value = 3
$ws.Range("W52").Value = 1
$ws.Range("X52").Value = 5
$ws.Range("Y52").Value = 5
$ws.Range("Z52").Value = 5
$ws.Range("AA52").Value = 5
$ws.Range("AB52").Value = 5
$ws.Range("A1").Copy()
$ws.Range("A52:AB52").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("AC52:AL52").PasteSpecial(-4122)

# Row 53
$ws.Range("A53").Value = 5285
$ws.Range("B53").Value = 3
$ws.Range("D53").Value = 2
$ws.Range("F53").Value = 5
$ws.Range("G53").Value = 5
$ws.Range("I53").Value = 4
$ws.Range("J53").Value = 4
$ws.Range("K53").Value = 3
$ws.Range("L53").Value = 2
$ws.Range("M53").Value = 2
$ws.Range("N53").Value = 4
$ws.Range("P53").Value = 2
$ws.Range("R53").Value = 4
$ws.Range("T53").Value = 4
$ws.Range("X53").Value = 4
$ws.Range("Y53").Value = 4
$ws.Range("Z53").Value = 4
$ws.Range("AB53").Value = 5
$ws.Range("AC53").Value = 5
$ws.Range("AD53").Value = 2
$ws.Range("AF53").Value = 3
$ws.Range("AG53").Value = 3
$ws.Range("AH53").Value = 4
$ws.Range("AI53").Value = 5
$ws.Range("A1").Copy()
$ws.Range("A53:B53").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G53").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I53").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("K53").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("M53").PasteSpecial(-4122)

# Row 54
$ws.Range("A54").Value = 1133
$ws.Range("B54").Value = 5
$ws.Range("C54").Value = 5
$ws.Range("F54").Value = 5
$ws.Range("H54").Value = 4
$ws.Range("J54").Value = 1
$ws.Range("L54").Value = 2
$ws.Range("M54").Value = 3
$ws.Range("N54").Value = 3
$ws.Range("O54").Value = 4
$ws.Range("P54").Value = 5
$ws.Range("Q54").Value = 5
$ws.Range("S54").Value = 3
$ws.Range("U54").Value = 4
$ws.Range("V54").Value = 4
$ws.Range("X54").Value = 4
$ws.Range("AA54").Value = 4
$ws.Range("AC54").Value = 3
$ws.Range("AE54").Value = 3
$ws.Range("AH54").Value = 3
$ws.Range("A1").Copy()
$ws.Range("A54:B54").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("O54").PasteSpecial(-4122)

# Row 55
$ws.Range("A55").Value = 4352
$ws.Range("B55").Value = 4
$ws.Range("C55").Value = 4
$ws.Range("G55").Value = 3
$ws.Range("I55").Value = 3
$ws.Range("K55").Value = 5
$ws.Range("M55").Value = 1
$ws.Range("N55").Value = 2
$ws.Range("O55").Value = 4
$ws.Range("P55").Value = 3
$ws.Range("Q55").Value = 5
$ws.Range("T55").Value = 4
$ws.Range("U55").Value = 3
$ws.Range("V55").Value = 5
$ws.Range("W55").Value = 3
$ws.Range("X55").Value = 1
$ws.Range("Y55").Value = 2
$ws.Range("AC55").Value = 4
$ws.Range("AD55").Value = 5
$ws.Range("AE55").Value = 3
$ws.Range("AF55").Value = 5
$ws.Range("AG55").Value = 3
$ws.Range("A1").Copy()
$ws.Range("A55:B55").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G55").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I55").PasteSpecial(-4122)

# Row 56
$ws.Range("A56").Value = 8096
$ws.Range("B56").Value = 5
$ws.Range("E56").Value = 5
$ws.Range("F56").Value = 2
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 4
$ws.Range("I56").Value = 4
$ws.Range("J56").Value = 5
$ws.Range("K56").Value = 4
$ws.Range("L56").Value = 3
$ws.Range("M56").Value = 3
$ws.Range("N56").Value = 4
$ws.Range("P56").Value = 2
$ws.Range("Q56").Value = 3
$ws.Range("R56").Value = 4
$ws.Range("S56").Value = 4
$ws.Range("W56").Value = 2
$ws.Range("X56").Value = 4
$ws.Range("Y56").Value = 3
$ws.Range("Z56").Value = 4
$ws.Range("AA56").Value = 4
$ws.Range("AB56").Value = 4
$ws.Range("AD56").Value = 4
$ws.Range("AE56").Value = 4
$ws.Range("AF56").Value = 3
$ws.Range("AG56").Value = 4
$ws.Range("AI56").Value = 5
$ws.Range("A1").Copy()
$ws.Range("A56:B56").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G56").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I56").PasteSpecial(-4122)

# Row 57
$ws.Range("A57").Value = 2290
$ws.Range("B57").Value = 3
$ws.Range("C57").Value = 4
$ws.Range("D57").Value = 5
$ws.Range("H57").Value = 4
$ws.Range("I57").Value = 3
$ws.Range("J57").Value = 4
$ws.Range("N57").Value = 4
$ws.Range("P57").Value = 3
$ws.Range("R57").Value = 1
$ws.Range("S57").Value = 4
$ws.Range("T57").Value = 5
$ws.Range("U57").Value = 3
$ws.Range("V57").Value = 5
$ws.Range("W57").Value = 2
$ws.Range("X57").Value = 4
$ws.Range("AA57").Value = 3
$ws.Range("AC57").Value = 5
$ws.Range("AE57").Value = 5
$ws.Range("AG57").Value = 5
$ws.Range("AH57").Value = 4
$ws.Range("A1").Copy()
$ws.Range("A57:B57").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I57").PasteSpecial(-4122)

# Row 58
$ws.Range("A58").Value = 9803
$ws.Range("B58").Value = 5
$ws.Range("C58").Value = 4
$ws.Range("D58").Value = 3
$ws.Range("E58").Value = 2
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 4
$ws.Range("I58").Value = 5
$ws.Range("J58").Value = 4
$ws.Range("K58").Value = 3
$ws.Range("M58").Value = 4
$ws.Range("N58").Value = 4
$ws.Range("O58").Value = 2
$ws.Range("R58").Value = 3
$ws.Range("S58").Value = 4
$ws.Range("T58").Value = 5
$ws.Range("U58").Value = 5
$ws.Range("V58").Value = 5
$ws.Range("W58").Value = 4
$ws.Range("Y58").Value = 4
$ws.Range("Z58").Value = 2
$ws.Range("AA58").Value = 3
$ws.Range("AB58").Value = 4
$ws.Range("AC58").Value = 4
$ws.Range("AF58").Value = 4
$ws.Range("AI58").Value = 3
$ws.Range("A1").Copy()
$ws.Range("A58:B58").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G58").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I58").PasteSpecial(-4122)

# Row 59
$ws.Range("A59").Value = 4445
$ws.Range("B59").Value = 4
$ws.Range("D59").Value = 5
$ws.Range("H59").Value = 4
$ws.Range("J59").Value = 4
$ws.Range("K59").Value = 3
$ws.Range("L59").Value = 4
$ws.Range("M59").Value = 3
$ws.Range("N59").Value = 3
$ws.Range("O59").Value = 3
$ws.Range("P59").Value = 4
$ws.Range("Q59").Value = 4
$ws.Range("R59").Value = 3
$ws.Range("S59").Value = 2
$ws.Range("U59").Value = 1
$ws.Range("V59").Value = 5
$ws.Range("W59").Value = 5
$ws.Range("X59").Value = 4
$ws.Range("Y59").Value = 4
$ws.Range("Z59").Value = 4
$ws.Range("AA59").Value = 4
$ws.Range("AB59").Value = 2
$ws.Range("AC59").Value = 1
$ws.Range("AD59").Value = 4
$ws.Range("AE59").Value = 4
$ws.Range("AF59").Value = 5
$ws.Range("AG59").Value = 5
$ws.Range("AH59").Value = 4
$ws.Range("AI59").Value = 3
$ws.Range("AJ59").Value = 1
$ws.Range("A1").Copy()
$ws.Range("A59:B59").PasteSpecial(-4122)

# Row 60
$ws.Range("A60").Value = 1176
$ws.Range("B60").Value = 4
$ws.Range("C60").Value = 4
$ws.Range("D60").Value = 5
$ws.Range("F60").Value = 3
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 4
$ws.Range("I60").Value = 4
$ws.Range("J60").Value = 3
$ws.Range("K60").Value = 4
$ws.Range("L60").Value = 5
$ws.Range("M60").Value = 3
$ws.Range("N60").Value = 4
$ws.Range("O60").Value = 3
$ws.Range("P60").Value = 4
$ws.Range("Q60").Value = 4
$ws.Range("R60").Value = 5
$ws.Range("S60").Value = 5
$ws.Range("T60").Value = 3
$ws.Range("V60").Value = 5
$ws.Range("W60").Value = 4
$ws.Range("X60").Value = 5
$ws.Range("Y60").Value = 4
$ws.Range("Z60").Value = 4
$ws.Range("AA60").Value = 2
$ws.Range("AB60").Value = 3
$ws.Range("AC60").Value = 5
$ws.Range("AE60").Value = 5
$ws.Range("AG60").Value = 3
$ws.Range("AI60").Value = 3
$ws.Range("AJ60").Value = 4
$ws.Range("A1").Copy()
$ws.Range("A60:B60").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G60").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I60").PasteSpecial(-4122)

# Set row height for the newly added rows to match source formatting (ht="19")
$ws.Range("A45:AL60").RowHeight = 19

# ---- Update view/selection state ----
# Select on PreProcessed first (not the final active tab)...
$ws3.Activate()
$ws3.Range("A19:B48").Select()

# ...then re-activate Data-Org and select there last, so it ends up the
# active tab (matches tabSelected on sheet1 in the target file).
$ws.Activate()
$ws.Range("A54").Select()
